# Edit script implementing the diff for Listado_Modulos2.xlsx
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Opción de Menú vs Funcion param"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# sheetView change (cosmetic)
$ws1.Application.ActiveWindow.ScrollRow = 300
$ws1.Range("D308").Select()

# Row 202: re-style whole row with highlighted ("fill2") style, values unchanged
$ws1.Range("A2").Copy()                 # style donor isn't ideal; replaced below
